# Cleaned up the image folder: removed the stimulus that used to live on
# row 18 (target/old/j, stimuli/img_a9acb.png) and replaced it with a new
# catch-trial row referencing stimuli/catch_18.jpg. The stimuli that used
# to occupy rows 19-23 each shift up into the row above them (their L..V
# columns move to row-1), and the old catch-trial row (23, formerly
# stimuli/catch_20_stairs.jpg) is reused to host the stimulus that used to
# be in row 22 (stimuli/img_u2o6z.png), now tagged as a regular "new" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("H18").ClearContents()
$ws.Range("I18").ClearContents()
$ws.Range("J18").Value = "catch"
$ws.Range("K18").Value = "f"
$ws.Range("L18").Value = "stimuli/catch_18.jpg"
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("O18").ClearContents()
$ws.Range("P18").ClearContents()
$ws.Range("Q18").ClearContents()
$ws.Range("R18").ClearContents()
$ws.Range("S18").ClearContents()
$ws.Range("T18").ClearContents()
$ws.Range("U18").ClearContents()
$ws.Range("V18").ClearContents()

# Row 19
$ws.Range("H19").Value = "bedrooms"
$ws.Range("I19").Value = "target"
$ws.Range("J19").Value = "old"
$ws.Range("K19").Value = "j"
$ws.Range("L19").Value = "stimuli/img_a9acb.png"
$ws.Range("M19").Value = 77.11428571428571
$ws.Range("N19").Value = 58.42857142857143
$ws.Range("O19").Value = 67.77142857142857
$ws.Range("P19").Value = 35
$ws.Range("Q19").Value = 7
$ws.Range("R19").Value = 7
$ws.Range("S19").Value = 7
$ws.Range("T19").Value = 7
$ws.Range("U19").Value = 7
$ws.Range("V19").Value = 7

# Row 20
$ws.Range("H20").Value = "bedrooms"
$ws.Range("I20").Value = "target"
$ws.Range("J20").Value = "old"
$ws.Range("K20").Value = "j"
$ws.Range("L20").Value = "stimuli/img_oou46.png"
$ws.Range("M20").Value = 75.70270270270271
$ws.Range("N20").Value = 54.86486486486486
$ws.Range("O20").Value = 65.28378378378379
$ws.Range("P20").Value = 37
$ws.Range("Q20").Value = 6
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 6
$ws.Range("T20").Value = 6
$ws.Range("U20").Value = 6
$ws.Range("V20").Value = 6

# Row 21
$ws.Range("H21").Value = "bedrooms"
$ws.Range("I21").Value = "target"
$ws.Range("J21").Value = "old"
$ws.Range("K21").Value = "j"
$ws.Range("L21").Value = "stimuli/img_wyctg.png"
$ws.Range("M21").Value = 33.44736842105263
$ws.Range("N21").Value = 11.39473684210526
$ws.Range("O21").Value = 22.42105263157895
$ws.Range("P21").Value = 38
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 1
$ws.Range("T21").Value = 1
$ws.Range("U21").Value = 1
$ws.Range("V21").Value = 1

# Row 22
$ws.Range("H22").Value = "bedrooms"
$ws.Range("J22").Value = "new"
$ws.Range("K22").Value = "f"
$ws.Range("L22").Value = "stimuli/img_wijef.png"
$ws.Range("M22").Value = 69.875
$ws.Range("N22").Value = 48.025
$ws.Range("O22").Value = 58.95
$ws.Range("P22").Value = 40
$ws.Range("Q22").Value = 5
$ws.Range("R22").Value = 5
$ws.Range("S22").Value = 5
$ws.Range("T22").Value = 5
$ws.Range("U22").Value = 5
$ws.Range("V22").Value = 5

# Row 23
$ws.Range("H23").Value = "bedrooms"
$ws.Range("J23").Value = "new"
$ws.Range("K23").Value = "f"
$ws.Range("L23").Value = "stimuli/img_u2o6z.png"
$ws.Range("M23").Value = 58.6
$ws.Range("N23").Value = 38.2
$ws.Range("O23").Value = 48.40000000000001
$ws.Range("P23").Value = 30
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = 3
$ws.Range("S23").Value = 3
$ws.Range("T23").Value = 3
$ws.Range("U23").Value = 3
$ws.Range("V23").Value = 3

